$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 144 (date 2020-08-08 / serial 44051) with corrected values
$ws.Range("B144").Value = 823
$ws.Range("C144").Value = 11
$ws.Range("D144").Value = 708
$ws.Range("E144").Value = 6
$ws.Range("F144").Value = 20
$ws.Range("G144").Value = 15
$ws.Range("H144").Value = 35
$ws.Range("I144").Value = 327
$ws.Range("J144").Value = 346
$ws.Range("K144").Value = 17
$ws.Range("L144").Value = 327
$ws.Range("M144").Value = 362
$ws.Range("N144").Value = 368

# Add new row 145 (date 2020-08-09 / serial 44052), copying the date format from A144
$ws.Range("A144").Copy()
$ws.Range("A145").PasteSpecial(-4122)
$ws.Range("A145").Value = 44052

$ws.Range("B145").Value = 823
$ws.Range("C145").Value = 15
$ws.Range("D145").Value = 724
$ws.Range("E145").Value = 5
$ws.Range("F145").Value = 21
$ws.Range("G145").Value = 15
$ws.Range("H145").Value = 30
$ws.Range("I145").Value = 336
$ws.Range("J145").Value = 358
$ws.Range("K145").Value = 17
$ws.Range("L145").Value = 336
$ws.Range("M145").Value = 366
$ws.Range("N145").Value = 371
